# Edit workbook outputs-HGR-r202-archive3/g__CAG-353.xlsx
# - Shift header labels in C1/D1/E1:
#     C1: max -> prediction
#     D1: prediction -> rejection-f
#     E1: rejection-f -> max
# - For each data row (2..26):
#     C column changes from the numeric "max" score to the species string
#       (same string already present in column D)
#     E column changes from the species string to a new numeric
#       "rejection-f" score

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# --- Data rows ---
$species = "s__CAG-353 sp900066885"

$rejectionValues = @{
  2  = 0.9567629966920078
  3  = 0.9474654544315221
  4  = 0.9434117746537035
  5  = 0.9413959015881465
  6  = 0.9458009503186926
  7  = 0.951345646495416
  8  = 0.9464662950434277
  9  = 0.9497162996387529
  10 = 0.9532136824614197
  11 = 0.9383665977282368
  12 = 0.9415332750327353
  13 = 0.9413707399575134
  14 = 0.9522152192398741
  15 = 0.9551580111710214
  16 = 0.9562034511884671
  17 = 0.9447607681122862
  18 = 0.938828719981086
  19 = 0.953757791190362
  20 = 0.9489500972298968
  21 = 0.9400640046247047
  22 = 0.9453249575615528
  23 = 0.9278537454736135
  24 = 0.9316883707293466
  25 = 0.9272112157150467
  26 = 0.9434576046073513
}

foreach ($row in 2..26) {
  $ws.Range("C$row").Value = $species
  $ws.Range("E$row").Value = $rejectionValues[$row]
}
